# "Generate Report for Handback"
# The localization status report is regenerated after a handback event:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this text lives in a shared string also used by the Overview sheet,
#      so updating it once updates Overview!E2/F2 and both language sheets'
#      Status column automatically).
#   - The "Latest Handback DateTime" for each locale is refreshed to the
#     time of this handback run.
#   - The stale "handback file is not latest" Error Detail is cleared now
#     that the handback is in sync.
#   - Column widths are refreshed to fit the new (wider/narrower) content.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# --- Status: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed for this handback run ---
$wsZhCn.Range("L2").Value = "2017-02-09 14:01:15"
$wsDeDe.Range("L2").Value = "2017-02-09 14:01:38"

# --- Error Detail cleared now that the handback is up to date ---
$wsZhCn.Range("R2").Value = ""
$wsDeDe.Range("R2").Value = ""

# --- Refresh column widths to fit the regenerated content ---
$wsOverview.Range("E:E").ColumnWidth = 29.14437166849777
$wsOverview.Range("F:F").ColumnWidth = 29.14437166849777

$wsZhCn.Range("C:C").ColumnWidth = 29.14437166849777
$wsZhCn.Range("R:R").ColumnWidth = 12.913719813028965

$wsDeDe.Range("C:C").ColumnWidth = 29.14437166849777
$wsDeDe.Range("R:R").ColumnWidth = 12.913719813028965
